# Refresh the "cryptos" price list (GitHub Actions scheduled update).
# For every data row, update the Price (column D) and Volume(1h) (column E)
# text cells with the latest scraped figures. Two coin pairs (RenderToken /
# VeChain at rows 39-40, and Elrond / Decentraland at rows 49-50) swapped
# rank order in this run, so their Coin/Link/Price/Volume cells are
# rewritten in the new row order as well.
#
# Column D holds numeric-looking values (e.g. "0.9971", "29.919.38") that
# must stay plain text exactly as scraped (some even use '.' as a
# thousands separator), so each Price cell is pinned to the Text number
# format before the value is written — otherwise Excel's normal
# type-inference would silently convert them to floating point numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.919.38"
$ws.Range("E2").Value = "  -1.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.874.30"
$ws.Range("E3").Value = "  -2.37%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9971"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.93"
$ws.Range("E5").Value = "  -5.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9975"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4963"
$ws.Range("E7").Value = "  -4.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.52"
$ws.Range("E8").Value = "  -3.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2913"
$ws.Range("E9").Value = "  -2.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06604"
$ws.Range("E10").Value = "  -3.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.872.88"
$ws.Range("E11").Value = "  -2.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.71"
$ws.Range("E12").Value = "  -4.93%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07234"
$ws.Range("E13").Value = "  -1.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6661"
$ws.Range("E14").Value = "  -3.89%  "
$ws.Range("E15").Value = "  -2.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.861"
$ws.Range("E16").Value = "  -1.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.878.44"
$ws.Range("E17").Value = "  -1.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007869"
$ws.Range("E18").Value = "  -3.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9974"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.75"
$ws.Range("E20").Value = "  -2.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.113.58"
$ws.Range("E21").Value = "  -2.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9961"
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.742"
$ws.Range("E23").Value = "  -2.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.607"
$ws.Range("E24").Value = "  -3.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.053"
$ws.Range("E25").Value = "  -2.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "149.21"
$ws.Range("E26").Value = "  +1.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "137.97"
$ws.Range("E27").Value = "  -1.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.97"
$ws.Range("E28").Value = "  -2.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.906"
$ws.Range("E29").Value = "  -5.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.382"
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.185"
$ws.Range("E31").Value = "  -3.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08777"
$ws.Range("E32").Value = "  -1.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.958"
$ws.Range("E33").Value = "  -2.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05031"
$ws.Range("E34").Value = "  -2.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7074"
$ws.Range("E35").Value = "  -2.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.103"
$ws.Range("E36").Value = "  -5.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.665"
$ws.Range("E37").Value = "  -0.72%  "
$ws.Range("E38").Value = "  -6.96%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.177"
$ws.Range("E39").Value = "  -7.00%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01730"
$ws.Range("E40").Value = "  +1.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9274"
$ws.Range("E41").Value = "  -5.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4264"
$ws.Range("E42").Value = "  -2.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.747"
$ws.Range("E43").Value = "  -7.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9963"
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.70"
$ws.Range("E45").Value = "  -4.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.403"
$ws.Range("E46").Value = "  -3.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1253"
$ws.Range("E47").Value = "  -2.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05648"
$ws.Range("E48").Value = "  -1.62%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "32.24"
$ws.Range("E49").Value = "  -3.36%  "
$ws.Range("B50").Value = "Decentraland"
$ws.Range("C50").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3750"
$ws.Range("E50").Value = "  -3.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.155"
$ws.Range("E51").Value = "  -4.97%  "
